$d = $word.ActiveDocument

$pairs = @(
    @("585÷7=83, 4", "614÷5=122, 4"),
    @("419÷8=52, 3", "127÷2=63, 1"),
    @("997÷7=142, 3", "681÷4=170, 1"),
    @("790÷4=197, 2", "130÷6=21, 4"),
    @("268÷3=89, 1", "773÷3=257, 2"),
    @("126÷5=25, 1", "699÷5=139, 4"),
    @("635÷6=105, 5", "669÷4=167, 1"),
    @("782÷7=111, 5", "818÷7=116, 6"),
    @("431÷9=47, 8", "899÷8=112, 3"),
    @("562÷3=187, 1", "288÷9=32, 0"),
    @("882÷3=294, 0", "951÷9=105, 6"),
    @("984÷8=123, 0", "774÷3=258, 0"),
    @("504÷5=100, 4", "461÷9=51, 2"),
    @("577÷5=115, 2", "478÷4=119, 2"),
    @("675÷7=96, 3", "574÷8=71, 6"),
    @("827÷8=103, 3", "935÷9=103, 8"),
    @("868÷3=289, 1", "181÷6=30, 1"),
    @("745÷2=372, 1", "849÷9=94, 3"),
    @("752÷5=150, 2", "169÷8=21, 1"),
    @("642÷8=80, 2", "523÷6=87, 1"),
    @("744÷6=124, 0", "986÷3=328, 2"),
    @("706÷2=353, 0", "707÷8=88, 3"),
    @("109÷2=54, 1", "761÷5=152, 1"),
    @("119÷7=17, 0", "238÷3=79, 1"),
    @("477÷9=53, 0", "224÷2=112, 0")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
